$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.191267490386963
$ws.Range("B1").Value = 2.297987937927246
$ws.Range("C1").Value = 6.628082752227783
$ws.Range("D1").Value = 2.325729608535767
$ws.Range("E1").Value = 1.189179658889771
